$d = $word.ActiveDocument

# Range covering the first three paragraphs (the ones the diff touches):
#   "Blabla", "Gulgul", "bunglesome" (+ the _GoBack bookmark).
# Paragraph 4 (empty, trailing) and the sectPr are left completely alone.
$p1 = $d.Paragraphs.Item(1)
$p3 = $d.Paragraphs.Item(3)
$target = $d.Range($p1.Range.Start, $p3.Range.End)

$body = '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Blabla</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' + `
        '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Gulgul</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' + `
        '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>B</w:t></w:r><w:r><w:t>unglesome</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' + `
        '<w:p><w:r><w:t>M</w:t></w:r><w:r><w:t>uddy</w:t></w:r></w:p>' + `
        '<w:p><w:r><w:t>F</w:t></w:r><w:r><w:t>labbergast</w:t></w:r></w:p>' + `
        '<w:p><w:r><w:t>A</w:t></w:r><w:r><w:t>stonish</w:t></w:r></w:p>' + `
        '<w:p><w:r><w:t>freegan</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
       '<pkg:xmlData>' + `
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
       '<w:body>' + $body + '</w:body></w:document>' + `
       '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)

Write-Output "done"
